$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 479294.2
$ws.Range("J17").Value = 479294.2
$ws.Range("L17").Value = 1437882.6
$ws.Range("N17").Value = -1438218.6
$ws.Range("H32").Value = 5292
$ws.Range("J32").Value = 5292
$ws.Range("L32").Value = 5292
$ws.Range("N32").Value = -5944
$ws.Range("H33").Value = 152.35
$ws.Range("I33").Value = 152.35
$ws.Range("K33").Value = 152.35
$ws.Range("M33").Value = 76.65000000000001
$ws.Range("H34").Value = 3676.5715
$ws.Range("I34").Value = 3676.5715
$ws.Range("K34").Value = 3676.5715
$ws.Range("M34").Value = -3473.5715
$ws.Range("H36").Value = 3676.5715
$ws.Range("I36").Value = 3676.5715
$ws.Range("K36").Value = 3676.5715
$ws.Range("M36").Value = -2961.5715
$ws.Range("H40").Value = 5555.9
$ws.Range("I40").Value = 8516
$ws.Range("J40").Value = 2595.8
$ws.Range("K40").Value = 8516
$ws.Range("L40").Value = 2595.8
$ws.Range("M40").Value = -8341
$ws.Range("N40").Value = -2945.8
$ws.Range("H69").Value = 9299.691999999999
$ws.Range("I69").Value = 9180.091
$ws.Range("J69").Value = 9957.5
$ws.Range("K69").Value = 27540.273
$ws.Range("L69").Value = 29872.5
$ws.Range("M69").Value = -26666.273
$ws.Range("N69").Value = -31620.5
$ws.Range("H72").Value = 9299.691999999999
$ws.Range("I72").Value = 9180.091
$ws.Range("J72").Value = 9957.5
$ws.Range("K72").Value = 82620.819
$ws.Range("L72").Value = 89617.5
$ws.Range("M72").Value = -78252.819
$ws.Range("N72").Value = -98353.5
$ws.Range("H137").Value = 2891695.8
$ws.Range("I137").Value = 4214558
$ws.Range("K137").Value = 12643674
$ws.Range("M137").Value = -12641124
$ws.Range("H138").Value = 3405.22
$ws.Range("I138").Value = 1557.5938
$ws.Range("J138").Value = 4274.6914
$ws.Range("K138").Value = 4672.7814
$ws.Range("L138").Value = 12824.0742
$ws.Range("M138").Value = 467.2186000000002
$ws.Range("N138").Value = -23104.0742

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10406.974
$ws.Range("I32").Value = 10045.423
$ws.Range("K32").Value = 10045.423
$ws.Range("M32").Value = -9758.423000000001
$ws.Range("H45").Value = 1827.8182
$ws.Range("I45").Value = 1327.2
$ws.Range("K45").Value = 1327.2
$ws.Range("M45").Value = -950.2
$ws.Range("H61").Value = 6152.6665
$ws.Range("I61").Value = 5052.4443
$ws.Range("J61").Value = 8353.111000000001
$ws.Range("K61").Value = 5052.4443
$ws.Range("L61").Value = 8353.111000000001
$ws.Range("M61").Value = -4840.4443
$ws.Range("N61").Value = -8777.111000000001
$ws.Range("H74").Value = 64529.344
$ws.Range("I74").Value = 66503.45
$ws.Range("J74").Value = 3332
$ws.Range("K74").Value = 66503.45
$ws.Range("L74").Value = 3332
$ws.Range("M74").Value = -65629.45
$ws.Range("N74").Value = -5080
$ws.Range("H77").Value = 64529.344
$ws.Range("I77").Value = 66503.45
$ws.Range("J77").Value = 3332
$ws.Range("K77").Value = 332517.25
$ws.Range("L77").Value = 16660
$ws.Range("M77").Value = -328149.25
$ws.Range("N77").Value = -25396
$ws.Range("H88").Value = 945.3611
$ws.Range("I88").Value = 985.8
$ws.Range("J88").Value = 916.4761999999999
$ws.Range("K88").Value = 985.8
$ws.Range("L88").Value = 916.4761999999999
$ws.Range("M88").Value = -579.8
$ws.Range("N88").Value = -1728.4762
$ws.Range("H91").Value = 945.3611
$ws.Range("I91").Value = 985.8
$ws.Range("J91").Value = 916.4761999999999
$ws.Range("K91").Value = 985.8
$ws.Range("L91").Value = 916.4761999999999
$ws.Range("M91").Value = 418.2
$ws.Range("N91").Value = -3724.4762
$ws.Range("H136").Value = 6152.6665
$ws.Range("I136").Value = 5052.4443
$ws.Range("J136").Value = 8353.111000000001
$ws.Range("K136").Value = 15157.3329
$ws.Range("L136").Value = 25059.333
$ws.Range("M136").Value = -12607.3329
$ws.Range("N136").Value = -30159.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 649.5
$ws.Range("I22").Value = 613
$ws.Range("K22").Value = 613
$ws.Range("M22").Value = -440
$ws.Range("H96").Value = 19928
$ws.Range("I96").Value = 19928
$ws.Range("K96").Value = 19928
$ws.Range("M96").Value = -17182
$ws.Range("H107").Value = 3018.5715
$ws.Range("I107").Value = 3130.5
$ws.Range("J107").Value = 2999.9167
$ws.Range("K107").Value = 3130.5
$ws.Range("L107").Value = 2999.9167
$ws.Range("M107").Value = -1210.5
$ws.Range("N107").Value = -6839.9167
$ws.Range("H134").Value = 1945.9524
$ws.Range("I134").Value = 1576.3889
$ws.Range("J134").Value = 4163.3335
$ws.Range("K134").Value = 4729.1667
$ws.Range("L134").Value = 12490.0005
$ws.Range("M134").Value = -2194.1667
$ws.Range("N134").Value = -17560.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 742.3333
$ws.Range("I35").Value = 774
$ws.Range("K35").Value = 774
$ws.Range("M35").Value = -480
$ws.Range("H94").Value = 1665.5
$ws.Range("I94").Value = 1851.6
$ws.Range("J94").Value = 1479.4
$ws.Range("K94").Value = 1851.6
$ws.Range("L94").Value = 1479.4
$ws.Range("M94").Value = -1400.6
$ws.Range("N94").Value = -2381.4
$ws.Range("H99").Value = 1114513
$ws.Range("I99").Value = 2003082.8
$ws.Range("J99").Value = 3800.75
$ws.Range("K99").Value = 2003082.8
$ws.Range("L99").Value = 3800.75
$ws.Range("M99").Value = -2001584.8
$ws.Range("N99").Value = -6796.75
$ws.Range("H126").Value = 1114513
$ws.Range("I126").Value = 2003082.8
$ws.Range("J126").Value = 3800.75
$ws.Range("K126").Value = 6009248.4
$ws.Range("L126").Value = 11402.25
$ws.Range("M126").Value = -6006778.4
$ws.Range("N126").Value = -16342.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100071780
$ws.Range("J37").Value = 100071780
$ws.Range("L37").Value = 300215340
$ws.Range("N37").Value = -300215564
$ws.Range("H39").Value = 1573.3529
$ws.Range("I39").Value = 853
$ws.Range("J39").Value = 2019.2858
$ws.Range("K39").Value = 2559
$ws.Range("L39").Value = 6057.857400000001
$ws.Range("M39").Value = -2265
$ws.Range("N39").Value = -6645.857400000001
$ws.Range("H69").Value = 2799.6667
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2799.6667
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 8399.000100000001
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -10021.0001
$ws.Range("H72").Value = 2799.6667
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2799.6667
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 25197.0003
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -33309.0003
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14617
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13674
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28752
$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 10000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -83760
$ws.Range("H98").Value = 282.2
$ws.Range("I98").Value = 282.2
$ws.Range("K98").Value = 846.5999999999999
$ws.Range("M98").Value = 651.4000000000001
$ws.Range("H137").Value = 2402.8
$ws.Range("I137").Value = 1886.75
$ws.Range("J137").Value = 2746.8333
$ws.Range("K137").Value = 5660.25
$ws.Range("L137").Value = 8240.499899999999
$ws.Range("M137").Value = -560.25
$ws.Range("N137").Value = -18440.4999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 18333
$ws.Range("I99").Value = 18333
$ws.Range("K99").Value = 18333
$ws.Range("M99").Value = -16087
$ws.Range("H126").Value = 7678.375
$ws.Range("I126").Value = 7813.5
$ws.Range("K126").Value = 23440.5
$ws.Range("M126").Value = -20970.5
$ws.Range("H132").Value = 24746.568
$ws.Range("J132").Value = 10935.546
$ws.Range("L132").Value = 32806.638
$ws.Range("N132").Value = -37866.638
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4124.3335
$ws.Range("I35").Value = 3949.2
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 3949.2
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -3613.2
$ws.Range("N35").Value = -5672
$ws.Range("H40").Value = 16242.125
$ws.Range("I40").Value = 16242.125
$ws.Range("K40").Value = 16242.125
$ws.Range("M40").Value = -16106.125
$ws.Range("H132").Value = 6623.025
$ws.Range("I132").Value = 5530.1665
$ws.Range("K132").Value = 16590.4995
$ws.Range("M132").Value = -14060.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2921.5518
$ws.Range("I122").Value = 2388.5217
$ws.Range("K122").Value = 7165.5651
$ws.Range("M122").Value = -4715.5651
$ws.Range("H126").Value = 1248.4117
$ws.Range("I126").Value = 1052.6
$ws.Range("K126").Value = 3157.8
$ws.Range("M126").Value = -687.7999999999997
$ws.Range("H132").Value = 1083.68
$ws.Range("I132").Value = 1092.7084
$ws.Range("K132").Value = 3278.1252
$ws.Range("M132").Value = -748.1251999999999
$ws.Range("H136").Value = 590648.4
$ws.Range("I136").Value = 627313.9399999999
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 1881941.82
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -1879391.82
$ws.Range("N136").Value = -17097
